# adding hp to charactor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: change B26 (percentage) from 0.5 (50%) to 1 (100%)
$ws.Range("B26").Value = 1

# Row 28: B28 gets a value of 1 (100%), C28 gets the new note "HP = 5"
# Match the existing percentage number format used by the other score cells in column B
$ws.Range("B28").NumberFormat = $ws.Range("B26").NumberFormat
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = "HP = 5"

# Update the selected/active cell shown in the sheet view to C29
$ws.Range("C29").Select()
